$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45001
$ws.Range("D4").Value = 45050
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 14000
$ws.Range("S4").Value = 778
$ws.Range("D5").Value = 45037
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 16000
$ws.Range("S5").Value = 889
$ws.Range("D6").Value = 45099
$ws.Range("M6").Value = 40
$ws.Range("N6").Value = 22000
$ws.Range("O6").Value = 22000
$ws.Range("P6").Value = 22000
$ws.Range("S6").Value = 1222
$ws.Range("D7").Value = 45036
$ws.Range("R7").Value = 'Región Metropolitana'
$ws.Range("D8").Value = 45028
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 18000
$ws.Range("O8").Value = 18000
$ws.Range("P8").Value = 18000
$ws.Range("S8").Value = 1000
$ws.Range("D9").Value = 45049
$ws.Range("M9").Value = 80
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 15000
$ws.Range("S9").Value = 833
$ws.Range("D10").Value = 45021
$ws.Range("R10").Value = 'Provincia de Los Andes'
$ws.Range("D11").Value = 45041
$ws.Range("D12").Value = 45030
$ws.Range("M12").Value = 40
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 18000
$ws.Range("S12").Value = 1000
$ws.Range("D13").Value = 45096
$ws.Range("M13").Value = 50
$ws.Range("N13").Value = 23000
$ws.Range("O13").Value = 23000
$ws.Range("P13").Value = 23000
$ws.Range("S13").Value = 1278
$ws.Range("D14").Value = 45062
$ws.Range("M14").Value = 90
$ws.Range("N14").Value = 13000
$ws.Range("O14").Value = 14000
$ws.Range("P14").Value = 13444
$ws.Range("S14").Value = 747
$ws.Range("D15").Value = 45014
$ws.Range("M15").Value = 30
$ws.Range("D16").Value = 45043
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range("S16").Value = 833
$ws.Range("D17").Value = 44999
$ws.Range("M17").Value = 60
$ws.Range("N17").Value = 17000
$ws.Range("O17").Value = 18000
$ws.Range("P17").Value = 17500
$ws.Range("Q17").Value = '$/caja 18 kilos'
$ws.Range("R17").Value = 'Región Metropolitana'
$ws.Range("S17").Value = 972
$ws.Range("T17").Value = 18
$ws.Range("D18").Value = 45020
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 15000
$ws.Range("Q18").Value = '$/caja 16 kilos'
$ws.Range("R18").Value = 'Provincia de Los Andes'
$ws.Range("S18").Value = 938
$ws.Range("T18").Value = 16
$ws.Range("D19").Value = 45033
$ws.Range("O19").Value = 16000
$ws.Range("P19").Value = 15500
$ws.Range("S19").Value = 861
$ws.Range("D21").Value = 45044
$ws.Range("M21").Value = 60
$ws.Range("N21").Value = 15000
$ws.Range("O21").Value = 15000
$ws.Range("P21").Value = 15000
$ws.Range("S21").Value = 833
$ws.Range("D22").Value = 45091
$ws.Range("M22").Value = 50
$ws.Range("N22").Value = 22000
$ws.Range("O22").Value = 22000
$ws.Range("P22").Value = 22000
$ws.Range("S22").Value = 1222